# Generate Report for Handoff
# Updates the localization-status report: marks statuses "Ready for handoff"
# (instead of "Handed back: in sync with en-US") and refreshes the
# handoff/generate timestamps, plus narrows the now-shorter status columns.

$wb = $excel.ActiveWorkbook

# Closest value achievable to the target column width (~17.216 chars) given
# this engine's 1/6-character snapping granularity on ColumnWidth.
$newColWidth = 16 + 1/3

# --- Sheet "Overview" ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-09-04 19:04:26"
$ws.Columns.Item(5).ColumnWidth = $newColWidth
$ws.Columns.Item(6).ColumnWidth = $newColWidth

# --- Sheet "zh-cn" ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-09-04 19:04:21"
$ws.Columns.Item(3).ColumnWidth = $newColWidth

# --- Sheet "de-de" ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-09-04 19:04:26"
$ws.Columns.Item(3).ColumnWidth = $newColWidth
